# "Generate Report for Handoff" — b.md has just been freshly handed off
# again (a new xliff pair was generated for zh-cn / de-de), so its status
# flips from "Handed back: in sync with en-US" to "Ready for handoff" on
# every sheet, and the per-language detail sheets record the new handoff
# file name/timestamp plus a note that the handback file is now stale.

$wb = $excel.ActiveWorkbook

$statusNew = "Ready for handoff"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/636eaff0e8357c328696c65c20378a3bfe7017f4/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/02280d8173e82bf41a3976fdc85318992e1e47ab/e2e/b.md."

# Helper: write a literal "True"/"False" piece of text into a cell without
# Excel's automatic boolean coercion kicking in (typing the bare word would
# turn the cell into a real boolean). Build it as a text formula, then
# collapse the formula down to its cached text value in place.
function Set-TextBoolean($range, [string]$text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the b.md file. Update its zh-cn/de-de status
# columns and the "Latest HO Xliff Generate Date" column.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusNew
$overview.Range("F3").Value = $statusNew
$overview.Range("G3").Value = "2016-08-24 06:38:08"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the b.md file.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusNew
Set-TextBoolean $zhcn.Range("F3") "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-24 06:37:57"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Range("P1").ColumnWidth = 39.16

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the b.md file.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusNew
Set-TextBoolean $dede.Range("F3") "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-24 06:38:08"
$dede.Range("P3").Value = $errorDetail
$dede.Range("P1").ColumnWidth = 39.16
